$p = $ppt.ActivePresentation

# Slide 1
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange
$titlePara = $tr.Paragraphs(1,1)
$titlePara.Text = "Juniper Green completes solar component of 133 MW Solapur hybrid power project - pv magazine India"
$titlePara.Font.Size = 28
$body1Para = $tr.Paragraphs(2,1)
$body1Para.Text = "Juniper Green has completed the solar component of the Solapur hybrid power project."
$body2Para = $tr.Paragraphs(3,1)
$body2Para.Text = "The total capacity of the Solapur hybrid power project is 133 MW."
# restore autofit-computed shape height to keep box size stable
$shape.Height = 324

# Slide 2
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange
$titlePara = $tr.Paragraphs(1,1)
$titlePara.Text = "Juniper Green completes solar component of 133 MW Solapur hybrid power project - pv magazine India"
$titlePara.Font.Size = 28
$body1Para = $tr.Paragraphs(2,1)
$body1Para.Text = "The project is located in Solapur, India."
$body2Para = $tr.Paragraphs(3,1)
$body2Para.Text = "The hybrid power project includes both solar and other energy components."
# restore autofit-computed shape height to keep box size stable
$shape.Height = 324

# Slide 3
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange
$titlePara = $tr.Paragraphs(1,1)
$titlePara.Text = "Juniper Green completes solar component of 133 MW Solapur hybrid power project - pv magazine India"
$titlePara.Font.Size = 28
$body1Para = $tr.Paragraphs(2,1)
$body1Para.Text = "The completion of the solar component contributes to renewable energy generation in the region."
$body2Para = $tr.Paragraphs(3,1)
$body2Para.Text = "Juniper Green is focused on expanding its renewable energy portfolio."
# restore autofit-computed shape height to keep box size stable
$shape.Height = 324

# Slide 4
$s = $p.Slides.Item(4)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange
$titlePara = $tr.Paragraphs(1,1)
$titlePara.Text = "Juniper Green completes solar component of 133 MW Solapur hybrid power project - pv magazine India"
$titlePara.Font.Size = 28
$body1Para = $tr.Paragraphs(2,1)
$body1Para.Text = "The project aligns with India's goals for increasing renewable energy capacity."
$body2Para = $tr.Paragraphs(3,1)
$body2Para.Text = "The hybrid power project is part of a broader initiative to enhance energy security."
# restore autofit-computed shape height to keep box size stable
$shape.Height = 324

# Slide 5
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange
$titlePara = $tr.Paragraphs(1,1)
$titlePara.Text = "Juniper Green completes solar component of 133 MW Solapur hybrid power project - pv magazine India"
$titlePara.Font.Size = 28
$body1Para = $tr.Paragraphs(2,1)
$body1Para.Text = "The solar component is a significant step in the development of the hybrid project."
$body2Para = $tr.Paragraphs(3,1)
$body2Para.Text = "Completion of this project may lead to further investments in renewable energy infrastructure."
# restore autofit-computed shape height to keep box size stable
$shape.Height = 324

# Slide 6
$s = $p.Slides.Item(6)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange
$titlePara = $tr.Paragraphs(1,1)
$titlePara.Text = "Juniper Green completes solar component of 133 MW Solapur hybrid power project - pv magazine India"
$titlePara.Font.Size = 28
$body1Para = $tr.Paragraphs(2,1)
$body1Para.Text = "The project showcases advancements in solar technology and hybrid systems."
$body2Para = $tr.Paragraphs(3,1)
$body2Para.Text = "Juniper Green's efforts contribute to the reduction of carbon emissions in the energy sector."
# restore autofit-computed shape height to keep box size stable
$shape.Height = 324
